# Add new worksheets "sheet2" and "sheet4", each seeded with a copy of the
# header row (rollno, firstname, lastname, nickname, marks) from "sheet1",
# plus a brand-new blank worksheet "sheet6" -- reusing data/files that
# already exist in the workbook as well as adding new ones.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# sheet2: new sheet placed right after sheet1, header row copied over
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet2.Name = "sheet2"
$sheet1.Range("A1:E1").Copy($sheet2.Range("A1"))

# sheet4: another new sheet, also seeded from the existing header row
$sheet4 = $wb.Worksheets.Add($null, $sheet2)
$sheet4.Name = "sheet4"
$sheet1.Range("A1:E1").Copy($sheet4.Range("A1"))

# sheet6: a fresh, empty sheet for new data
$sheet6 = $wb.Worksheets.Add($null, $sheet4)
$sheet6.Name = "sheet6"

# keep the original sheet active/selected, as it was before the edit
$sheet1.Activate()
$sheet1.Select()
